# Weekly update: insert the newest week's price records at the top of the
# data (rows 3 and 4), pushing all previously existing rows down by two
# positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (old rows 3..21) down by two rows by inserting
# two fresh rows above the current row 3.
$ws.Rows("3:4").Insert()

# Static values shared by every data row in this sheet.
$mercado    = "Agrícola del Norte S.A. de Arica"
$region     = "Arica y Parinacota"
$codreg     = 15
$catId      = 100112006
$categoria  = "Repollo"
$variedad   = "Copenhague"
$unidadCom  = "$/unidad"
$origen     = "Región de Arica y Parinacota"
$kgOUnidad  = 1
$clasif     = "Hortaliza"

# New row 3: Calidad "Segunda"
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = $mercado
$ws.Range("C3").Value = $region
$ws.Range("D3").Value = 45282
$ws.Range("E3").Value = $codreg
$ws.Range("F3").Value = $catId
$ws.Range("G3").Value = $categoria
$ws.Range("H3").Value = $variedad
$ws.Range("I3").Value = "Segunda"
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 650
$ws.Range("L3").Value = 700
$ws.Range("M3").Value = 672
$ws.Range("N3").Value = $unidadCom
$ws.Range("O3").Value = $origen
$ws.Range("P3").Value = 672
$ws.Range("Q3").Value = $kgOUnidad
$ws.Range("R3").Value = $clasif

# New row 4: Calidad "Tercera"
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = $mercado
$ws.Range("C4").Value = $region
$ws.Range("D4").Value = 45282
$ws.Range("E4").Value = $codreg
$ws.Range("F4").Value = $catId
$ws.Range("G4").Value = $categoria
$ws.Range("H4").Value = $variedad
$ws.Range("I4").Value = "Tercera"
$ws.Range("J4").Value = 850
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 600
$ws.Range("M4").Value = 571
$ws.Range("N4").Value = $unidadCom
$ws.Range("O4").Value = $origen
$ws.Range("P4").Value = 571
$ws.Range("Q4").Value = $kgOUnidad
$ws.Range("R4").Value = $clasif

Write-Host "Inserted two new rows and populated data."
